# Auto-generated script to apply 2025-09-16 violent crime data updates
# across the 'violent-crime-full-year.xlsx' workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 4778
$ws.Range('L3').Value = 5152
$ws.Range('C4').Value = 1870
$ws.Range('H4').Value = 1761
$ws.Range('J4').Value = 1874
$ws.Range('L4').Value = 1265
$ws.Range('L5').Value = 304
$ws.Range('K6').Value = 9115
$ws.Range('L6').Value = 4371
$ws.Range('C7').Value = 28414
$ws.Range('H7').Value = 26077
$ws.Range('J7').Value = 29350
$ws.Range('K7').Value = 27571
$ws.Range('L7').Value = 15870

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 308
$ws.Range('L3').Value = 355
$ws.Range('L6').Value = 276
$ws.Range('L7').Value = 1058

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L2').Value = 108
$ws.Range('L7').Value = 352

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('L2').Value = 85
$ws.Range('L7').Value = 218

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L2').Value = 173
$ws.Range('L5').Value = 18
$ws.Range('L6').Value = 163
$ws.Range('L7').Value = 588

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L3').Value = 98
$ws.Range('L4').Value = 14
$ws.Range('L7').Value = 308

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L2').Value = 132
$ws.Range('L7').Value = 520
$ws.Range('L8').Value = 1058
$ws.Range('L17').Value = 30
$ws.Range('L19').Value = 436
$ws.Range('L20').Value = 397
$ws.Range('L25').Value = 93
$ws.Range('L26').Value = 16
$ws.Range('L27').Value = 141
$ws.Range('L29').Value = 870
$ws.Range('L31').Value = 159
$ws.Range('L37').Value = 588
$ws.Range('L40').Value = 43
$ws.Range('L42').Value = 517
$ws.Range('L47').Value = 110
$ws.Range('L48').Value = 205
$ws.Range('L50').Value = 78
$ws.Range('L51').Value = 197
$ws.Range('L52').Value = 319
$ws.Range('L54').Value = 332
$ws.Range('C63').Value = 295
$ws.Range('H63').Value = 312
$ws.Range('J63').Value = 226
$ws.Range('K63').Value = 167
$ws.Range('L63').Value = 46
$ws.Range('L65').Value = 308
$ws.Range('L67').Value = 547
$ws.Range('L79').Value = 418
$ws.Range('L80').Value = 50
$ws.Range('L83').Value = 352
$ws.Range('L85').Value = 814
$ws.Range('L87').Value = 47
$ws.Range('L89').Value = 229
$ws.Range('L92').Value = 44
$ws.Range('L95').Value = 218
$ws.Range('L96').Value = 176
$ws.Range('C101').Value = 28414
$ws.Range('H101').Value = 26077
$ws.Range('J101').Value = 29350
$ws.Range('K101').Value = 27571
$ws.Range('L101').Value = 15870

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('L2').Value = 62
$ws.Range('L7').Value = 159

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L3').Value = 213
$ws.Range('L6').Value = 125
$ws.Range('L7').Value = 547

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L6').Value = 162
$ws.Range('L7').Value = 332

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L2').Value = 260
$ws.Range('L6').Value = 229
$ws.Range('L7').Value = 870

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('L3').Value = 51
$ws.Range('L7').Value = 205

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L2').Value = 154
$ws.Range('L3').Value = 136
$ws.Range('L7').Value = 436

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L2').Value = 148
$ws.Range('L4').Value = 39
$ws.Range('L7').Value = 517

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L6').Value = 51
$ws.Range('L7').Value = 176

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L2').Value = 137
$ws.Range('L7').Value = 418

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L2').Value = 122
$ws.Range('L7').Value = 397

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range('L6').Value = 10
$ws.Range('L7').Value = 30

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L3').Value = 173
$ws.Range('L7').Value = 520

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('L5').Value = 2
$ws.Range('L7').Value = 93

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('L6').Value = 25
$ws.Range('L7').Value = 110

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('L3').Value = 19
$ws.Range('L7').Value = 78

$ws = $wb.Worksheets.Item('East Village')
$ws.Range('L3').Value = 2
$ws.Range('L7').Value = 16

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('L3').Value = 43
$ws.Range('L7').Value = 132

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('L2').Value = 17
$ws.Range('L7').Value = 44

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('L3').Value = 68
$ws.Range('L6').Value = 61
$ws.Range('L7').Value = 229

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('L2').Value = 40
$ws.Range('L7').Value = 141

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('L6').Value = 46
$ws.Range('L7').Value = 197

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L3').Value = 330
$ws.Range('L7').Value = 814

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('L6').Value = 26
$ws.Range('L7').Value = 50

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range('L4').Value = 2
$ws.Range('L7').Value = 43

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('L3').Value = 102
$ws.Range('L7').Value = 319

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('L6').Value = 17
$ws.Range('L7').Value = 47

